# Auto-generated edit script: update market-price derived cells per commit diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1045.1666
$ws.Range("J17").Value = 1045.1666
$ws.Range("L17").Value = 3135.4998
$ws.Range("N17").Value = -3471.4998
$ws.Range("H58").Value = 1562.3334
$ws.Range("I58").Value = 344
$ws.Range("J58").Value = 3999
$ws.Range("K58").Value = 1032
$ws.Range("L58").Value = 11997
$ws.Range("M58").Value = -882
$ws.Range("N58").Value = -12297
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("H98").Value = 1279.0834
$ws.Range("J98").Value = 3495
$ws.Range("L98").Value = 3495
$ws.Range("N98").Value = -6491
$ws.Range("H122").Value = 1279.0834
$ws.Range("J122").Value = 3495
$ws.Range("L122").Value = 10485
$ws.Range("N122").Value = -15385
$ws.Range("H141").Value = 6469.143
$ws.Range("I141").Value = 4852.1113
$ws.Range("J141").Value = 9379.799999999999
$ws.Range("K141").Value = 14556.3339
$ws.Range("L141").Value = 28139.4
$ws.Range("M141").Value = -9376.333899999998
$ws.Range("N141").Value = -38499.39999999999
$ws.Range("M86").ClearContents()
$ws.Range("M89").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8985.109
$ws.Range("I32").Value = 3604.709
$ws.Range("K32").Value = 3604.709
$ws.Range("M32").Value = -3317.709
$ws.Range("H47").Value = 36540.363
$ws.Range("I47").Value = 36994
$ws.Range("J47").Value = 36439.555
$ws.Range("K47").Value = 36994
$ws.Range("L47").Value = 36439.555
$ws.Range("N47").Value = -37889.555
$ws.Range("M47").Value = -36269

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 33816.582
$ws.Range("J4").Value = 47542.715
$ws.Range("L4").Value = 47542.715
$ws.Range("N4").Value = -47766.715
$ws.Range("H22").Value = 518.2143
$ws.Range("I22").Value = 518.9
$ws.Range("K22").Value = 518.9
$ws.Range("M22").Value = -168.9
$ws.Range("H31").Value = 5340.607
$ws.Range("I31").Value = 2706
$ws.Range("J31").Value = 8380.538
$ws.Range("K31").Value = 2706
$ws.Range("L31").Value = 8380.538
$ws.Range("M31").Value = -2411
$ws.Range("N31").Value = -8970.538
$ws.Range("H34").Value = 5340.607
$ws.Range("I34").Value = 2706
$ws.Range("J34").Value = 8380.538
$ws.Range("K34").Value = 2706
$ws.Range("L34").Value = 8380.538
$ws.Range("M34").Value = -2504
$ws.Range("N34").Value = -8784.538
$ws.Range("H62").Value = 17876.363
$ws.Range("I62").Value = 16821.334
$ws.Range("J62").Value = 18606.77
$ws.Range("K62").Value = 16821.334
$ws.Range("L62").Value = 18606.77
$ws.Range("M62").Value = -16197.334
$ws.Range("N62").Value = -19854.77
$ws.Range("H65").Value = 17876.363
$ws.Range("I65").Value = 16821.334
$ws.Range("J65").Value = 18606.77
$ws.Range("K65").Value = 84106.67
$ws.Range("L65").Value = 93033.85000000001
$ws.Range("M65").Value = -80986.67
$ws.Range("N65").Value = -99273.85000000001
$ws.Range("H86").Value = 2749.9412
$ws.Range("J86").Value = 3144.8333
$ws.Range("L86").Value = 3144.8333
$ws.Range("N86").Value = -5390.8333
$ws.Range("H89").Value = 2749.9412
$ws.Range("J89").Value = 3144.8333
$ws.Range("L89").Value = 15724.1665
$ws.Range("N89").Value = -26956.1665
$ws.Range("H107").Value = 1773.6666
$ws.Range("I107").Value = 994
$ws.Range("J107").Value = 3333
$ws.Range("K107").Value = 994
$ws.Range("L107").Value = 3333
$ws.Range("M107").Value = 926
$ws.Range("N107").Value = -7173
$ws.Range("H132").Value = 1586.5
$ws.Range("I132").Value = 1308.375
$ws.Range("K132").Value = 3925.125
$ws.Range("M132").Value = -1395.125
$ws.Range("H134").Value = 1743.871
$ws.Range("I134").Value = 1743.871
$ws.Range("K134").Value = 5231.613
$ws.Range("M134").Value = -2696.613

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H118").Value = 2000
$ws.Range("I118").Value = 2000
$ws.Range("K118").Value = 6000
$ws.Range("M118").Value = -4757

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 568.3333
$ws.Range("I97").Value = 511.7143
$ws.Range("K97").Value = 511.7143
$ws.Range("M97").Value = -15.71429999999998
$ws.Range("H102").Value = 2062.1482
$ws.Range("I102").Value = 1270.381
$ws.Range("K102").Value = 1270.381
$ws.Range("M102").Value = 351.6189999999999
$ws.Range("H132").Value = 1867
$ws.Range("I132").Value = 1867
$ws.Range("K132").Value = 5601
$ws.Range("M132").Value = -3071

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 897.9
$ws.Range("I46").Value = 887.2
$ws.Range("J46").Value = 908.6
$ws.Range("K46").Value = 887.2
$ws.Range("L46").Value = 908.6
$ws.Range("M46").Value = -699.2
$ws.Range("N46").Value = -1284.6
$ws.Range("H100").Value = 2405.6072
$ws.Range("I100").Value = 2069.5454
$ws.Range("J100").Value = 2623.0588
$ws.Range("K100").Value = 2069.5454
$ws.Range("L100").Value = 2623.0588
$ws.Range("M100").Value = -1528.5454
$ws.Range("N100").Value = -3705.0588
$ws.Range("H124").Value = 33999.668
$ws.Range("J124").Value = 33999.668
$ws.Range("L124").Value = 33999.668
$ws.Range("N124").Value = -43819.668
$ws.Range("H130").Value = 52332.332
$ws.Range("J130").Value = 52332.332
$ws.Range("L130").Value = 52332.332
$ws.Range("N130").Value = -62372.332

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 13999.5
$ws.Range("I41").Value = 13999.5
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 13999.5
$ws.Range("L41").Value = 0
$ws.Range("H105").Value = 49999.5
$ws.Range("J105").Value = 49999.5
$ws.Range("L105").Value = 49999.5
$ws.Range("N105").Value = -56987.5
$ws.Range("H107").Value = 1455.8085
$ws.Range("J107").Value = 1237.6428
$ws.Range("L107").Value = 3712.9284
$ws.Range("N107").Value = -7552.928400000001
$ws.Range("H126").Value = 2110.6206
$ws.Range("I126").Value = 2193.6
$ws.Range("J126").Value = 1592
$ws.Range("K126").Value = 6580.799999999999
$ws.Range("L126").Value = 4776
$ws.Range("M126").Value = -4110.799999999999
$ws.Range("N126").Value = -9716
$ws.Range("H132").Value = 2312.5833
$ws.Range("I132").Value = 2312.5833
$ws.Range("K132").Value = 6937.749899999999
$ws.Range("M132").Value = -4407.749899999999
$ws.Range("M41").Value = -13609.5
$ws.Range("N41").ClearContents()

